$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 411 (the "「諸行無常」..." entry) and shift all following rows up by one.
$ws.Rows.Item(411).Delete()
